# fix(publipostage): Correct status name
#
# - "bleu" -> "noir" (statut_label, column B)
# - statut_name wording (column C) updated from
#   "résultat et / ou publication posté..." to
#   "résultat postés ou publiés..." (and matching "pas de résultat..." variant)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = $ws.Cells

$cells.Replace("bleu", "noir", -4163, 1, $false, $false)

$cells.Replace("résultat et / ou publication posté dans les 36 mois", "résultat postés ou publiés dans les 36 mois", -4163, 1, $false, $false)
$cells.Replace("résultat et / ou publication posté dans les 12 mois", "résultat postés ou publiés dans les 12 mois", -4163, 1, $false, $false)
$cells.Replace("résultat et / ou publication posté", "résultat postés ou publiés", -4163, 1, $false, $false)
$cells.Replace("pas de résultat ni de publication", "pas de résultat postés ni publiés", -4163, 1, $false, $false)
